# Regenerate merged AHB files
# - rename the "_old" / "_new" header-row suffixes to "_FV2410" / "_FV2504"
# - turn the data range into a real Excel Table (ListObject)
# - freeze the header row

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- 1. header row rename -------------------------------------------------
# Columns A:J describe the "old" (FV2410) side, K is the literal "diff"
# marker column, and L:U describe the "new" (FV2504) side.
$headersFV2410 = @(
    "Segmentname",
    "Segmentgruppe",
    "Segment",
    "Datenelement",
    "Segment ID",
    "Code",
    "Qualifier",
    "Beschreibung",
    "Bedingungsausdruck",
    "Bedingung"
)

for ($i = 0; $i -lt $headersFV2410.Length; $i++) {
    $ws.Cells.Item(1, $i + 1).Value = ($headersFV2410[$i] + "_FV2410")
}

# Column K ("diff") is unchanged.

for ($i = 0; $i -lt $headersFV2410.Length; $i++) {
    $ws.Cells.Item(1, $i + 12).Value = ($headersFV2410[$i] + "_FV2504")
}

# --- 2. turn the range into a table ---------------------------------------
$lo = $ws.ListObjects.Add(1, $ws.Range("A1:U73"), $null, 1)
$lo.Name = "Table1"

# --- 3. freeze the header row ----------------------------------------------
$ws.Range("A2").Select()
$excel.ActiveWindow.FreezePanes = $true
